$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.209.68'
$ws.Range('E2').Value = '  -0.81%  '
$ws.Range('D3').Value = '2.597.56'
$ws.Range('E3').Value = '  -0.81%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'592.56"
$ws.Range('E5').Value = '  -1.88%  '
$ws.Range('D6').Value = "'149.76"
$ws.Range('E6').Value = '  -3.31%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = "'0.541"
$ws.Range('E8').Value = '  -1.43%  '
$ws.Range('D9').Value = '2.596.60'
$ws.Range('E9').Value = '  -0.78%  '
$ws.Range('D10').Value = "'0.134"
$ws.Range('E10').Value = '  +7.11%  '
$ws.Range('E11').Value = '  -0.66%  '
$ws.Range('D12').Value = "'5.18"
$ws.Range('E12').Value = '  -1.40%  '
$ws.Range('E13').Value = '  -2.61%  '
$ws.Range('D14').Value = "'27.14"
$ws.Range('E14').Value = '  -3.30%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').Value = "'0.0000184"
$ws.Range('E15').Value = '  -0.83%  '
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = '3.074.56'
$ws.Range('E16').Value = '  -0.75%  '
$ws.Range('D17').Value = '67.077.38'
$ws.Range('E17').Value = '  -0.89%  '
$ws.Range('D18').Value = '2.598.80'
$ws.Range('E18').Value = '  -0.57%  '
$ws.Range('D19').Value = "'368.45"
$ws.Range('E19').Value = '  -0.62%  '
$ws.Range('D20').Value = "'10.99"
$ws.Range('E20').Value = '  -2.41%  '
$ws.Range('D21').Value = "'7.34"
$ws.Range('E21').Value = '  -3.23%  '
$ws.Range('D22').Value = "'4.18"
$ws.Range('E22').Value = '  -3.01%  '
$ws.Range('D23').Value = "'4.74"
$ws.Range('E23').Value = '  -5.01%  '
$ws.Range('D24').Value = "'2.01"
$ws.Range('E24').Value = '  -4.56%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Value = "'1.00"
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('B26').Value = 'Aptos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D26').Value = "'9.83"
$ws.Range('E26').Value = '  -1.88%  '
$ws.Range('B27').Value = 'Litecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D27').Value = "'66.17"
$ws.Range('E27').Value = '  -0.44%  '
$ws.Range('D28').Value = '2.731.55'
$ws.Range('E28').Value = '  -0.74%  '
$ws.Range('E29').Value = '  +0.35%  '
$ws.Range('D30').Value = "'580.10"
$ws.Range('E30').Value = '  -0.70%  '
$ws.Range('D31').Value = '0.0₃0989'
$ws.Range('E31').Value = '  -5.70%  '
$ws.Range('E32').Value = '  -5.67%  '
$ws.Range('D33').Value = "'7.65"
$ws.Range('E33').Value = '  -3.53%  '
$ws.Range('E34').Value = '  -2.69%  '
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').Value = "'0.125"
$ws.Range('E36').Value = '  -4.54%  '
$ws.Range('E37').Value = '  -2.06%  '
$ws.Range('D38').Value = "'157.44"
$ws.Range('E38').Value = '  +0.53%  '
$ws.Range('D39').Value = "'19.02"
$ws.Range('E39').Value = '  -2.46%  '
$ws.Range('E40').Value = '  -1.85%  '
$ws.Range('D41').Value = "'1.85"
$ws.Range('E41').Value = '  -0.26%  '
$ws.Range('D42').Value = "'5.18"
$ws.Range('E42').Value = '  -3.43%  '
$ws.Range('E43').Value = '  -3.82%  '
$ws.Range('D44').Value = "'17.08"
$ws.Range('E44').Value = '  +3.88%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value = "'40.12"
$ws.Range('E46').Value = '  -2.49%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = "'152.42"
$ws.Range('E47').Value = '  -2.53%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0281'
$ws.Range('E48').Value = '  -2.67%  '
$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').Value = "'3.63"
$ws.Range('E49').Value = '  -3.33%  '
$ws.Range('B50').Value = 'Optimism'
$ws.Range('C50').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D50').Value = "'1.67"
$ws.Range('E50').Value = '  -4.86%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = "'0.0774"
$ws.Range('E51').Value = '  -1.64%  '
